$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 18

$ws.Cells.Item($row, 1).Value = 131113915            # A - Id
$ws.Cells.Item($row, 2).Value = 57881                # B - Taxonsorteringsordning
$ws.Cells.Item($row, 4).Value = "NT"                 # D - Rödlistade
$ws.Cells.Item($row, 5).Value = 100049               # E - TaxonId
$ws.Cells.Item($row, 6).Value = "Spillkråka"          # F - Artnamn
$ws.Cells.Item($row, 7).Value = "Dryocopus martius"  # G - Vetenskapligt namn
$ws.Cells.Item($row, 8).Value = "(Linnaeus, 1758)"   # H - Auktor
$ws.Cells.Item($row, 9).Value = "'1"                 # I - Antal (kept as text)
$ws.Cells.Item($row, 11).Value = "adult"             # K - Ålder-Stadium
$ws.Cells.Item($row, 14).Value = "observerad"        # N - Metod
$ws.Cells.Item($row, 16).Value = "Trollknabbarna, Dockasberg, Nb"  # P - Lokalnamn
$ws.Cells.Item($row, 17).Value = 816315              # Q - Ost
$ws.Cells.Item($row, 18).Value = 7375643             # R - Nord
$ws.Cells.Item($row, 19).Value = 25                  # S - Noggrannhet
$ws.Cells.Item($row, 20).Value = "Norrbotten"        # T - Län
$ws.Cells.Item($row, 21).Value = "Överkalix"          # U - Kommun
$ws.Cells.Item($row, 22).Value = "Norrbotten"        # V - Provins
$ws.Cells.Item($row, 23).Value = "Överkalix"          # W - Socken
$ws.Cells.Item($row, 25).Value = "'2025-07-31"       # Y - Startdatum (text)
$ws.Cells.Item($row, 26).Value = "'10:02"            # Z - Starttid (text)
$ws.Cells.Item($row, 27).Value = "'2025-07-31"       # AA - Slutdatum (text)
$ws.Cells.Item($row, 28).Value = "'10:05"            # AB - Sluttid (text)
$ws.Cells.Item($row, 30).Value = $false              # AD - Ej återfunnen
$ws.Cells.Item($row, 31).Value = $false              # AE - Osäker artbestämning
$ws.Cells.Item($row, 33).Value = $false              # AG - Ospontan
$ws.Cells.Item($row, 49).Value = "Markus  Kristoffersson"  # AW - Rapportör
$ws.Cells.Item($row, 50).Value = "Markus  Kristoffersson"  # AX - Observatörer
